# Fruta / hortaliza, semanal
# Rows 2-13 have their D/H/I/J/K/L/M/P field-sets permuted among rows
# (the identifying columns A,B,C,E,F,G,N,O,Q,R stay fixed per row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record" of values for each row.
$cols = @("D", "H", "I", "J", "K", "L", "M", "P")

# Snapshot the current (before) values of the moving columns for every row.
# Value2 is used (instead of Value) so dates come back as raw numeric
# serials rather than DateTime/formatted-text, matching the stored XML.
$before = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Row -> row the new values come from (i.e. after[r] = before[perm[r]]).
$perm = @{
    2  = 6
    3  = 10
    4  = 5
    5  = 3
    6  = 13
    7  = 4
    8  = 8
    9  = 12
    10 = 7
    11 = 9
    12 = 11
    13 = 2
}

foreach ($r in $perm.Keys) {
    $src = $perm[$r]
    $srcVals = $before[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
